$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer data (x, y, z) to be written into rows 2-31.
# This replaces the existing rows 2-21 data: the first two original
# data rows are dropped, the remaining rows shift up by two, and the
# sheet is extended with new rows recorded on May 9th.
$data = @(
    @(-0.3062458038330078, 0.5210974216461182, -0.6826221942901611),
    @(-0.0336570739746093, -0.3159241676330566, -2.89306116104126),
    @(2.157181739807129, 0.2950803935527801, -5.19840145111084),
    @(0.2004642486572265, -0.8322474360466003, -0.1898527145385742),
    @(-2.500746965408325, 0.4417282938957214, 0.8889427185058594),
    @(-3.396578788757324, 1.112873911857605, 1.645373344421387),
    @(-4.896120071411133, 3.986298561096192, 8.319293975830078),
    @(-3.247961044311523, 1.173810482025146, -5.774084091186523),
    @(-0.8262515068054199, -2.918009757995605, 2.994272947311401),
    @(-2.458094835281372, -0.3888258635997772, -0.6840519905090332),
    @(-1.984735012054444, 2.097556114196777, -3.859614849090576),
    @(-10.1377477645874, 4.743382453918457, -1.631471395492554),
    @(3.672384262084961, -0.4976945519447326, 9.970148086547852),
    @(-3.353246688842773, 0.9116206765174866, -0.0151057243347167),
    @(-1.351092576980591, 1.431697487831116, 2.426383018493652),
    @(-4.564205169677734, 2.659675121307373, 9.756929397583008),
    @(-0.1535825729370117, -2.909769535064697, -4.346651077270508),
    @(-3.01861572265625, -0.8810358047485352, -0.5351743698120117),
    @(-3.367526054382324, 0.915987193584442, -1.626443386077881),
    @(-1.118759155273438, 2.792432069778442, -4.963344097137451),
    @(-1.557756900787354, -0.5582034587860107, -0.2619988918304443),
    @(1.715949058532715, -1.576748490333557, 5.096891403198242),
    @(-0.8243503570556641, 0.5943599939346313, 1.927432060241699),
    @(-0.2981023788452148, 1.024843096733093, 0.8517363667488098),
    @(0.4371089935302734, 0.3337190449237823, -0.154114544391632),
    @(-0.5641984939575195, -0.3292053341865539, -0.326197862625122),
    @(0.131052017211914, 0.5107872486114502, 0.0540084838867187),
    @(-0.1777238845825195, 0.4102384448051452, 0.1352127194404602),
    @(0.2971744537353515, 0.0503720641136169, -0.09675005078315729),
    @(-0.1711950302124023, -0.0577961653470993, 0.0454300940036773)
)

$rowCount = $data.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
